$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 3112.0667
$ws.Range("I9").Value = 246.3
$ws.Range("J9").Value = 8843.6
$ws.Range("K9").Value = 246.3
$ws.Range("L9").Value = 8843.6
$ws.Range("M9").Value = -77.30000000000001
$ws.Range("N9").Value = -9181.6

$ws.Range("H11").Value = 31011.5
$ws.Range("I11").Value = 31011.5
$ws.Range("K11").Value = 31011.5
$ws.Range("M11").Value = -30871.5

$ws.Range("H86").Value = 90911550
$ws.Range("I86").Value = 166668670
$ws.Range("J86").Value = 3018.8
$ws.Range("K86").Value = 166668670
$ws.Range("L86").Value = 3018.8
$ws.Range("M86").Value = -166667547
$ws.Range("N86").Value = -5264.8

$ws.Range("H89").Value = 90911550
$ws.Range("I89").Value = 166668670
$ws.Range("J89").Value = 3018.8
$ws.Range("K89").Value = 833343350
$ws.Range("L89").Value = 15094
$ws.Range("M89").Value = -833337734
$ws.Range("N89").Value = -26326

$ws.Range("H100").Value = 103335590
$ws.Range("J100").Value = 166669520
$ws.Range("L100").Value = 166669520
$ws.Range("N100").Value = -166670602

$ws.Range("H132").Value = 2998
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws.Range("H137").Value = 4349694.5
$ws.Range("J137").Value = 5001876
$ws.Range("L137").Value = 15005628
$ws.Range("N137").Value = -15010728

$ws.Range("H141").Value = 1123.2142
$ws.Range("I141").Value = 776.5
$ws.Range("J141").Value = 1990
$ws.Range("K141").Value = 2329.5
$ws.Range("L141").Value = 5970
$ws.Range("M141").Value = 2850.5
$ws.Range("N141").Value = -16330

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 18940.666
$ws.Range("I2").Value = 26748.5
$ws.Range("J2").Value = 3325
$ws.Range("K2").Value = 26748.5
$ws.Range("L2").Value = 3325
$ws.Range("M2").Value = -26635.5
$ws.Range("N2").Value = -3551

$ws.Range("H45").Value = 70575
$ws.Range("J45").Value = 4065
$ws.Range("L45").Value = 4065
$ws.Range("N45").Value = -4819

$ws.Range("H61").Value = 2033730.9
$ws.Range("I61").Value = 5111.2085
$ws.Range("K61").Value = 5111.2085
$ws.Range("M61").Value = -4899.2085

$ws.Range("H74").Value = 1115266
$ws.Range("I74").Value = 2123.111
$ws.Range("K74").Value = 2123.111
$ws.Range("M74").Value = -1249.111

$ws.Range("H77").Value = 1115266
$ws.Range("I77").Value = 2123.111
$ws.Range("K77").Value = 10615.555
$ws.Range("M77").Value = -6247.555

$ws.Range("H101").Value = 59000
$ws.Range("J101").Value = 59000
$ws.Range("L101").Value = 59000
$ws.Range("N101").Value = -65490

$ws.Range("H105").Value = 124999
$ws.Range("J105").Value = 124999
$ws.Range("L105").Value = 124999
$ws.Range("N105").Value = -131987

$ws.Range("H116").Value = 18940.666
$ws.Range("I116").Value = 26748.5
$ws.Range("J116").Value = 3325
$ws.Range("K116").Value = 26748.5
$ws.Range("L116").Value = 3325
$ws.Range("M116").Value = -24454.5
$ws.Range("N116").Value = -7913

$ws.Range("H119").Value = 73398.8
$ws.Range("J119").Value = 73398.8
$ws.Range("L119").Value = 73398.8
$ws.Range("N119").Value = -83074.8

$ws.Range("H136").Value = 2033730.9
$ws.Range("I136").Value = 5111.2085
$ws.Range("K136").Value = 15333.6255
$ws.Range("M136").Value = -12783.6255

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 18940.666
$ws.Range("I3").Value = 26748.5
$ws.Range("J3").Value = 3325
$ws.Range("K3").Value = 26748.5
$ws.Range("L3").Value = 3325
$ws.Range("M3").Value = -26634.5
$ws.Range("N3").Value = -3553

$ws.Range("H40").Value = 60448
$ws.Range("J40").Value = 60448
$ws.Range("L40").Value = 60448
$ws.Range("N40").Value = -60978

$ws.Range("H57").Value = 47500
$ws.Range("J57").Value = 47500
$ws.Range("L57").Value = 47500
$ws.Range("N57").Value = -48940

$ws.Range("H64").Value = 1868
$ws.Range("I64").Value = 2475
$ws.Range("J64").Value = 1694.5714
$ws.Range("K64").Value = 2475
$ws.Range("L64").Value = 1694.5714
$ws.Range("M64").Value = -2250
$ws.Range("N64").Value = -2144.5714

$ws.Range("H67").Value = 1868
$ws.Range("I67").Value = 2475
$ws.Range("J67").Value = 1694.5714
$ws.Range("K67").Value = 2475
$ws.Range("L67").Value = 1694.5714
$ws.Range("M67").Value = -1695
$ws.Range("N67").Value = -3254.5714

$ws.Range("H99").Value = 4941.7095
$ws.Range("I99").Value = 7156.5625
$ws.Range("J99").Value = 2579.2
$ws.Range("K99").Value = 7156.5625
$ws.Range("L99").Value = 2579.2
$ws.Range("M99").Value = -5658.5625
$ws.Range("N99").Value = -5575.2

$ws.Range("H134").Value = 42860460
$ws.Range("I134").Value = 2726.6667
$ws.Range("K134").Value = 8180.000100000001
$ws.Range("M134").Value = -5645.000100000001

$ws.Range("H136").Value = 47500
$ws.Range("J136").Value = 47500
$ws.Range("L136").Value = 47500
$ws.Range("N136").Value = -57700

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3079.25
$ws.Range("I31").Value = 2706.3333
$ws.Range("J31").Value = 3120.6853
$ws.Range("K31").Value = 2706.3333
$ws.Range("L31").Value = 3120.6853
$ws.Range("M31").Value = -2411.3333
$ws.Range("N31").Value = -3710.6853

$ws.Range("H34").Value = 3079.25
$ws.Range("I34").Value = 2706.3333
$ws.Range("J34").Value = 3120.6853
$ws.Range("K34").Value = 2706.3333
$ws.Range("L34").Value = 3120.6853
$ws.Range("M34").Value = -2504.3333
$ws.Range("N34").Value = -3524.6853

$ws.Range("H58").Value = 2901.923
$ws.Range("I58").Value = 2535
$ws.Range("J58").Value = 3268.8462
$ws.Range("K58").Value = 2535
$ws.Range("L58").Value = 3268.8462
$ws.Range("M58").Value = -2332
$ws.Range("N58").Value = -3674.8462

$ws.Range("H86").Value = 30739.416
$ws.Range("I86").Value = 16425
$ws.Range("J86").Value = 50779.6
$ws.Range("K86").Value = 16425
$ws.Range("L86").Value = 50779.6
$ws.Range("M86").Value = -15302
$ws.Range("N86").Value = -53025.6

$ws.Range("H89").Value = 30739.416
$ws.Range("I89").Value = 16425
$ws.Range("J89").Value = 50779.6
$ws.Range("K89").Value = 82125
$ws.Range("L89").Value = 253898
$ws.Range("M89").Value = -76509
$ws.Range("N89").Value = -265130

$ws.Range("H131").Value = 187000
$ws.Range("J131").Value = 187000
$ws.Range("L131").Value = 187000
$ws.Range("N131").Value = -197080

$ws.Range("H132").Value = 12348214
$ws.Range("I132").Value = 2024
$ws.Range("K132").Value = 6072
$ws.Range("M132").Value = -3542

$ws.Range("H136").Value = 2901.923
$ws.Range("I136").Value = 2535
$ws.Range("J136").Value = 3268.8462
$ws.Range("K136").Value = 7605
$ws.Range("L136").Value = 9806.5386
$ws.Range("M136").Value = -5055
$ws.Range("N136").Value = -14906.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 681456.4
$ws.Range("I4").Value = 838241.9
$ws.Range("K4").Value = 2514725.7
$ws.Range("M4").Value = -2514613.7

$ws.Range("H14").Value = 1819.5714
$ws.Range("I14").Value = 1819.5714
$ws.Range("K14").Value = 5458.7142
$ws.Range("M14").Value = -5285.7142

$ws.Range("H137").Value = 4250.4736
$ws.Range("I137").Value = 2576.6
$ws.Range("J137").Value = 6110.3335
$ws.Range("K137").Value = 7729.799999999999
$ws.Range("L137").Value = 18331.0005
$ws.Range("M137").Value = -2629.799999999999
$ws.Range("N137").Value = -28531.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 48000

$ws.Range("H65").Value = 48000

$ws.Range("H105").Value = 70670
$ws.Range("J105").Value = 70670
$ws.Range("L105").Value = 70670
$ws.Range("N105").Value = -77658

$ws.Range("H132").Value = 12531573
$ws.Range("J132").Value = 34458884
$ws.Range("L132").Value = 103376652
$ws.Range("N132").Value = -103381712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 10000
$ws.Range("I48").Value = 10000
$ws.Range("K48").Value = 10000
$ws.Range("M48").Value = -9339

$ws.Range("H68").Value = 3450.6667

$ws.Range("H71").Value = 3450.6667

$ws.Range("H100").Value = 2057.2856
$ws.Range("I100").Value = 1780.4
$ws.Range("K100").Value = 1780.4
$ws.Range("M100").Value = -1239.4

$ws.Range("H105").Value = 36947.2
$ws.Range("J105").Value = 36947.2
$ws.Range("L105").Value = 36947.2
$ws.Range("N105").Value = -43935.2

$ws.Range("H122").Value = 5446.4136
$ws.Range("I122").Value = 4000.182
$ws.Range("J122").Value = 6330.222
$ws.Range("K122").Value = 12000.546
$ws.Range("L122").Value = 18990.666
$ws.Range("M122").Value = -9550.545999999998
$ws.Range("N122").Value = -23890.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 29907.5
$ws.Range("J55").Value = 29907.5
$ws.Range("L55").Value = 29907.5
$ws.Range("N55").Value = -30461.5

$ws.Range("H82").Value = 26666.666
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 26666.666
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H107").Value = 125000580
$ws.Range("I107").Value = 764.4
$ws.Range("J107").Value = 333333600
$ws.Range("K107").Value = 2293.2
$ws.Range("L107").Value = 1000000800
$ws.Range("M107").Value = -373.1999999999998
$ws.Range("N107").Value = -1000004640
